$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws3 = $wb.Worksheets.Item("本地生活")
$ws4 = $wb.Worksheets.Item("全部类型")

$ws1.Range("F2").Value = 1709
$ws1.Range("F4").Value = 9916
$ws1.Range("F6").Value = 292
$ws1.Range("F10").Value = 1462
$ws1.Range("F13").Value = 1549
$ws1.Range("F15").Value = 342
$ws1.Range("F17").Value = 172
$ws1.Range("F18").Value = 453
$ws1.Range("F19").Value = 1138
$ws1.Range("F20").Value = 118
$ws1.Range("F23").Value = 67
$ws1.Range("F24").Value = 321
$ws1.Range("F26").Value = 292
$ws1.Range("F27").Value = 87
$ws1.Range("F29").Value = 664
$ws1.Range("F32").Value = 202
$ws1.Range("F34").Value = 211
$ws1.Range("F38").Value = 456
$ws1.Range("F39").Value = 686
$ws1.Range("F42").Value = 783
$ws1.Range("F43").Value = 355
$ws1.Range("F44").Value = 310
$ws1.Range("F45").Value = 344
$ws1.Range("F46").Value = 71
$ws1.Range("F47").Value = 341
$ws2.Range("F10").Value = 3
$ws2.Range("F12").Value = 54
$ws2.Range("F18").Value = 1045
$ws2.Range("F20").Value = 437
$ws2.Range("F22").Value = 310
$ws2.Range("F24").Value = 57
$ws2.Range("F27").Value = 9
$ws2.Range("F28").Value = 348
$ws2.Range("F34").Value = 147
$ws2.Range("F35").Value = 178
$ws2.Range("F36").Value = 40
$ws2.Range("F39").Value = 121
$ws2.Range("F41").Value = 54
$ws3.Range("F5").Value = 175
$ws3.Range("F6").Value = 2472
$ws3.Range("F7").Value = 3936
$ws3.Range("F8").Value = 41
$ws3.Range("F10").Value = 210
$ws3.Range("F11").Value = 159
$ws4.Range("F2").Value = 1709
$ws4.Range("F4").Value = 9916
$ws4.Range("F6").Value = 3936
$ws4.Range("F8").Value = 210
$ws4.Range("F9").Value = 210
$ws4.Range("F11").Value = 1549
$ws4.Range("F13").Value = 342
$ws4.Range("F15").Value = 3
$ws4.Range("F16").Value = 453
$ws4.Range("F17").Value = 1138
$ws4.Range("F18").Value = 118
$ws4.Range("F19").Value = 54
$ws4.Range("F22").Value = 67
$ws4.Range("F23").Value = 1045
$ws4.Range("F24").Value = 321
$ws4.Range("F27").Value = 292
$ws4.Range("F29").Value = 664
$ws4.Range("F31").Value = 57
$ws4.Range("F33").Value = 202
$ws4.Range("F34").Value = 348
$ws4.Range("F39").Value = 456
$ws4.Range("F41").Value = 686
$ws4.Range("F43").Value = 783
$ws4.Range("F44").Value = 178
$ws4.Range("F45").Value = 355
$ws4.Range("F46").Value = 40
$ws4.Range("F47").Value = 310
$ws4.Range("F48").Value = 344
$ws4.Range("F49").Value = 341
$ws4.Range("F50").Value = 54
